$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-05-08 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-05-09 Friday", 2) | Out-Null
$d.Content.Find.Execute("14×58=", $true, $false, $false, $false, $false, $true, 1, $false, "41×30=", 2) | Out-Null
$d.Content.Find.Execute("44×62=", $true, $false, $false, $false, $false, $true, 1, $false, "77×99=", 2) | Out-Null
$d.Content.Find.Execute("34×18=", $true, $false, $false, $false, $false, $true, 1, $false, "94×15=", 2) | Out-Null
$d.Content.Find.Execute("39×28=", $true, $false, $false, $false, $false, $true, 1, $false, "19×60=", 2) | Out-Null
$d.Content.Find.Execute("95×41=", $true, $false, $false, $false, $false, $true, 1, $false, "98×75=", 2) | Out-Null
$d.Content.Find.Execute("30×73=", $true, $false, $false, $false, $false, $true, 1, $false, "71×85=", 2) | Out-Null
$d.Content.Find.Execute("25×66=", $true, $false, $false, $false, $false, $true, 1, $false, "45×81=", 2) | Out-Null
$d.Content.Find.Execute("46×80=", $true, $false, $false, $false, $false, $true, 1, $false, "62×78=", 2) | Out-Null
$d.Content.Find.Execute("83×17=", $true, $false, $false, $false, $false, $true, 1, $false, "48×52=", 2) | Out-Null
$d.Content.Find.Execute("53×21=", $true, $false, $false, $false, $false, $true, 1, $false, "50×64=", 2) | Out-Null
$d.Content.Find.Execute("43×33=", $true, $false, $false, $false, $false, $true, 1, $false, "51×31=", 2) | Out-Null
$d.Content.Find.Execute("78×68=", $true, $false, $false, $false, $false, $true, 1, $false, "34×53=", 2) | Out-Null
$d.Content.Find.Execute("65×27=", $true, $false, $false, $false, $false, $true, 1, $false, "83×15=", 2) | Out-Null
$d.Content.Find.Execute("52×60=", $true, $false, $false, $false, $false, $true, 1, $false, "81×52=", 2) | Out-Null
$d.Content.Find.Execute("89×54=", $true, $false, $false, $false, $false, $true, 1, $false, "78×72=", 2) | Out-Null
$d.Content.Find.Execute("89×17=", $true, $false, $false, $false, $false, $true, 1, $false, "16×98=", 2) | Out-Null
$d.Content.Find.Execute("76×20=", $true, $false, $false, $false, $false, $true, 1, $false, "14×58=", 2) | Out-Null
$d.Content.Find.Execute("64×21=", $true, $false, $false, $false, $false, $true, 1, $false, "80×30=", 2) | Out-Null
$d.Content.Find.Execute("44×30=", $true, $false, $false, $false, $false, $true, 1, $false, "73×35=", 2) | Out-Null
$d.Content.Find.Execute("25×45=", $true, $false, $false, $false, $false, $true, 1, $false, "75×32=", 2) | Out-Null
$d.Content.Find.Execute("48×72=", $true, $false, $false, $false, $false, $true, 1, $false, "98×57=", 2) | Out-Null
$d.Content.Find.Execute("41×69=", $true, $false, $false, $false, $false, $true, 1, $false, "99×46=", 2) | Out-Null
$d.Content.Find.Execute("43×56=", $true, $false, $false, $false, $false, $true, 1, $false, "52×59=", 2) | Out-Null
$d.Content.Find.Execute("44×26=", $true, $false, $false, $false, $false, $true, 1, $false, "74×48=", 2) | Out-Null
$d.Content.Find.Execute("83×65=", $true, $false, $false, $false, $false, $true, 1, $false, "32×32=", 2) | Out-Null
